$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.946.11'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '2.432.90'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.158'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.97%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.329'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '67.850.84'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000169'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '333.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").Value = '0.0₃0809'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '419.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.105'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.295'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("E38").Value = '  +1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '129.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0709'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.478'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.556'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0914'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.42%  '
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '0.0₆0205'
$ws.Range("E50").Value = '  +5.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0429'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.07%  '
